$d = $word.ActiveDocument

# --- Hunk 1: the _GoBack bookmark that sat right before "BIG IDEAS FOR SMALL
#     BUSINESS" is gone in the new revision. ---
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Hunk 2: the run of spaces right before "THE BUSINESS MATCHMAKER" shrinks
#     from 39 to 37 characters (36 + 1, split into two runs), and the
#     _GoBack bookmark now sits between those two runs, right before
#     "THE BUSINESS MATCHMAKER". ---
$rng = $d.Content
$null = $rng.Find.Execute("THE BUSINESS MATCHMAKER", $true, $false, $false, `
                           $false, $false, $true, 1, $false, "", 0)
$target = $rng.Start
$endOfPhrase = $rng.End

# Pin both edges of the edit with scratch bookmarks first so the engine only
# rebuilds the runs strictly between them, leaving every other run (and its
# w:rsidR) untouched.
$d.Bookmarks.Add("TMP_EDIT_START", $d.Range($target - 2, $target - 2))
$d.Bookmarks.Add("TMP_EDIT_END", $d.Range($endOfPhrase, $endOfPhrase))

$d.Range($target - 2, $target).Delete()

$d.Bookmarks.Item("TMP_EDIT_START").Delete()
$d.Bookmarks.Item("TMP_EDIT_END").Delete()

# Re-locate "THE BUSINESS MATCHMAKER" now that two characters were removed.
$rng2 = $d.Content
$null = $rng2.Find.Execute("THE BUSINESS MATCHMAKER", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)
$target2 = $rng2.Start

# Re-add _GoBack one character before the phrase: this splits the remaining
# 37-space run into a 36-space run and a 1-space run, with the bookmark in
# between, exactly like the new revision.
$d.Bookmarks.Add("_GoBack", $d.Range($target2 - 1, $target2 - 1))
